# [Timesheet Calculator v2] Fix the bug that computes the credited
# minutes (col E, "Credited Regular Log [480 = 1 day]") and the excess
# minutes (col F, "Minutes in excess of 480; Sat/Sun Duties") during
# creation of the employee spreadsheet. Fixed testdata to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (16-Jun-2021): credited minutes and excess minutes were swapped.
# 480 rendered minutes on a weekday -> all 480 are credited, 0 excess.
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 480

# Row 22 (22-Jun-2021, Saturday duty): 120 rendered minutes on a
# Sat/Sun day should not count as "credited regular log" minutes; they
# are all excess minutes instead.
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 120

# C8 and E22 previously carried a stray one-off font/style variant
# that isn't used anywhere else in the sheet; normalize them back to
# the plain style shared by the rest of the sheet.
$ws.Range("C8").Font.Name = "Arial"
$ws.Range("E22").Font.Name = "Arial"

# Restore the cursor / selection to where the edit was made.
$ws.Range("F22").Select()
